$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the current row 14 (Biodiesel) to hold the
# new Municipal Solid Waste subfuel entries, shifting Biodiesel, Electricity
# and Heat down by two rows.
$ws.Rows.Item(14).Resize(2).Insert()

$ws.Cells.Item(14, 1).Value = "Municipal Solid Waste (Renewable)"
$ws.Cells.Item(14, 2).Value = "16_others"
$ws.Cells.Item(14, 3).Value = "16_03_municipal_solid_waste_renewable"

$ws.Cells.Item(15, 1).Value = "Municipal Solid Waste (Non-renewable)"
$ws.Cells.Item(15, 2).Value = "16_others"
$ws.Cells.Item(15, 3).Value = "16_04_municipal_solid_waste_nonrenewable"

# Adjust column widths to fit the new, longer content
$ws.Columns.Item(1).ColumnWidth = 37.28515625
$ws.Columns.Item(3).ColumnWidth = 42.42578125

# Update the selected cell to match the target state
$ws.Range("C21").Select()

$wb.Save()
